# Feb 19th offsets update (discobisco) - rename/restructure header rows on
# "Jersey Vitals" and "Jersey Colors" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Jersey Vitals" -----------------------------------------
# Shrinks from A1:AA1 (27 cols) down to A1:Y1 (25 cols): B/C swap places,
# several headers get renamed to their UPPER_SNAKE_CASE offset-name form,
# and the old "Uniform Address" / "Uniform File" columns (Y/Z) are dropped
# so the final UNIQUEID column shifts from AA1 to Y1.
$wsVitals = $wb.Worksheets.Item("Jersey Vitals")

$vitalsHeaders = @(
    "COLOR_LUMINANCE_LEVEL",
    "EDITION_CRC32B",
    "EDITION_NAME",
    "FILENAME",
    "HEADBAND_LOGO_TYPE",
    "IS_ALTERNATE",
    "IS_HOME",
    "IS_TEAM_CREATE_UNIFORM",
    "JACKET_WARMUP_CRC32B",
    "LOGO_BRAND",
    "LOGO_TYPE",
    "MYTEAM_INCLUDE",
    "NUMBER_ON_SHORTS",
    "OVERRIDE_NBA_PATCH",
    "SOCKS_AWAY",
    "SOCKS_COLOR_N#1",
    "SOCKS_COLOR_N#2",
    "SOCKS_COLOR_N#3",
    "SOCKS_COLOR_N#4",
    "SOCKS_COLOR_N#5",
    "SOCKS_HOME",
    "SPONSOR_PATCH",
    "TEAM",
    "TYPE",
    "UNIQUEID"
)

for ($i = 0; $i -lt $vitalsHeaders.Length; $i++) {
    $wsVitals.Cells.Item(1, $i + 1).Value = $vitalsHeaders[$i]
}

# Clear the two trailing columns (old Z1 "Uniform File" and AA1 "UNIQUEID")
# that are no longer part of the header row.
$wsVitals.Range("Z1:AA1").ClearContents()

# --- Sheet 2: "Jersey Colors" ------------------------------------------
# Grows from A1:W1 (23 cols) to A1:AA1 (27 cols): headers are renamed to
# UPPER_SNAKE_CASE offset names, the "Jersey Colors - X" prefix is dropped,
# and several "* Sock" / "* Accessory" columns are split out into their own
# "#SOCK" / "#ACCESSORY" suffixed columns.
$wsColors = $wb.Worksheets.Item("Jersey Colors")

$colorsHeaders = @(
    "ARM_ACCESSORY_COLOR_AWAY",
    "ARM_ACCESSORY_COLOR_HOME",
    "HEADBAND_COLOR_AWAY",
    "HEADBAND_COLOR_HOME",
    "LEG_ACCESSORY_COLOR_AWAY",
    "LEG_ACCESSORY_COLOR_HOME",
    "PRIMARY_COLOR",
    "PRIMARY_COLOR#ACCESSORY",
    "PRIMARY_COLOR#SOCK",
    "QUATERNARY_COLOR",
    "QUATERNARY_COLOR#ACCESSORY",
    "QUATERNARY_COLOR#SOCK",
    "QUINTARY_COLOR",
    "QUINTARY_COLOR#ACCESSORY",
    "QUINTARY_COLOR#SOCK",
    "SECONDARY_COLOR#ACCESSORY",
    "SECONDARY_COLOR#SOCK",
    "SECONDARY_COLOR_SHOES_COLOR",
    "SEXTARY_COLOR",
    "SHOE_PRIMARY_COLOR_AWAY",
    "SOCK_COLOR_AWAY",
    "SOCK_COLOR_HOME",
    "TERTIARY_COLOR",
    "TERTIARY_COLOR#ACCESSORY",
    "TERTIARY_COLOR#SOCK",
    "UNDER_SHIRT_COLOR_AWAY",
    "UNDER_SHIRT_COLOR_HOME"
)

for ($i = 0; $i -lt $colorsHeaders.Length; $i++) {
    $wsColors.Cells.Item(1, $i + 1).Value = $colorsHeaders[$i]
}
